$d = $word.ActiveDocument

# The edits below use Range.InsertXML so that run-splits, the moved
# "_GoBack" bookmark, and the lastRenderedPageBreak marker land exactly
# where the target OOXML expects them. Paragraphs are touched from the
# bottom of the document upward so that earlier Paragraphs.Item(n)
# indices stay valid even though paragraph 13 below gets split into
# three paragraphs.

# Paragraph 15: "A resource loader javascript class ..." -- reword,
# add the "This resource loader is provided by mtwilson-core-html5."
# run, and move the _GoBack bookmark to the end of this paragraph.
$p = $d.Paragraphs.Item(15)
$p.Range.InsertXML('<w:p><w:r><w:t xml:space="preserve">A resource loader javascript class in the UI is in charge of loading HTML, CSS, and JS files, ensuring that each file is loaded only once, and ensuring that a Javascript snippet that declares a dependency on external Javascript files is not run until all those files have been loaded. </w:t></w:r><w:r><w:t>This resource loader is provided by mtwilson-core-html5.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>') | Out-Null

# Paragraph 14: "Resource Loader" heading -- add lastRenderedPageBreak.
$p = $d.Paragraphs.Item(14)
$p.Range.InsertXML('<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Resource Loader</w:t></w:r></w:p>') | Out-Null

# Paragraph 13: "Any .jar file can add ..." -- split into three
# paragraphs (intro sentence, "First, ..." and "Second, ...").
$p = $d.Paragraphs.Item(13)
$p.Range.InsertXML('<w:p><w:r><w:t xml:space="preserve">There are two methods for providing static HTML5 content to the UI. </w:t></w:r></w:p><w:p><w:r><w:t>First, a</w:t></w:r><w:r><w:t xml:space="preserve">ny .jar file can add static HTML5 content to the UI by placing it under a “publicResources” folder. In the source tree this would be under the project’s src/main/resources/publicResources. A JAX-RS class </w:t></w:r><w:r><w:t xml:space="preserve">in mtwilson-core-html5 </w:t></w:r><w:r><w:t>provides read access to all files under “publicResources” in the classpath.</w:t></w:r><w:r><w:t xml:space="preserve"> However, the UI has to already “know” to load those resources from a specific URI relative to the application URI. </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Second, any feature can add static HTML5 content to the UI by placing it under a src/main/html5 folder in the Maven project layout and inheriting from mtwilson-core-feature-zip. See blueprint for mtwilson-core-html5 for more details. This content can then be discovered by and automatically integrated into the UI by including extensions for known extension points. </w:t></w:r></w:p>') | Out-Null

# Paragraph 11: "The kms-html5 project contains ..." -- reword to
# mtwilson-core-html5 phrasing.
$p = $d.Paragraphs.Item(11)
$p.Range.InsertXML('<w:p><w:r><w:t xml:space="preserve">The </w:t></w:r><w:r><w:t>mtwilson-core</w:t></w:r><w:r><w:t xml:space="preserve">-html5 project contains the static HTML5 resources for the main application. Other features are loaded and inserted </w:t></w:r><w:r><w:t>to the UI by</w:t></w:r><w:r><w:t xml:space="preserve"> this code.</w:t></w:r></w:p>') | Out-Null

# Paragraph 9: "The files for the splash screen are stored under ..."
# -- replace with the much shorter mtwilson-core-html5 sentence.
$p = $d.Paragraphs.Item(9)
$p.Range.InsertXML('<w:p><w:r><w:t xml:space="preserve">The splash screen is provided by mtwilson-core-html5. </w:t></w:r></w:p>') | Out-Null

# Paragraph 8: "The splash screen's style ..." -- "it can contain"
# becomes "it contains" (split across three runs).
$p = $d.Paragraphs.Item(8)
$p.Range.InsertXML('<w:p><w:r><w:t>The splash screen’s style should be consistent with the application style and it contain</w:t></w:r><w:r><w:t>s</w:t></w:r><w:r><w:t xml:space="preserve"> some elements of the application main page for continuity.</w:t></w:r></w:p>') | Out-Null

# Paragraph 2: "HTML5" title -- drop the _GoBack bookmark (it now
# lives at the end of the Resource Loader paragraph instead).
$p = $d.Paragraphs.Item(2)
$p.Range.InsertXML('<w:p><w:pPr><w:pStyle w:val="Title"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>HTML5</w:t></w:r></w:p>') | Out-Null

Write-Host "Paragraphs now:" $d.Paragraphs.Count
